$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Neg_Change")

$ws1.Range("A2").Value = "NESTLEIND"
$ws1.Range("B2").Value = 1263.9
$ws1.Range("C2").Value = 1311.6
$ws1.Range("D2").Value = 1261.3
$ws1.Range("E2").Value = 1289
$ws1.Range("F2").Value = 3171952
$ws1.Range("G2").Value = 6927367
$ws1.Range("H2").Value = -0.5421128980173853
$ws1.Range("I2").Value = "NESTLEIND"

$ws1.Range("A3").Value = "HDFCLIFE"
$ws1.Range("B3").Value = 740.75
$ws1.Range("C3").Value = 749.5
$ws1.Range("D3").Value = 736.05
$ws1.Range("E3").Value = 745.5
$ws1.Range("F3").Value = 4205888
$ws1.Range("G3").Value = 10209104
$ws1.Range("H3").Value = -0.5880257464318123
$ws1.Range("I3").Value = "HDFCLIFE"

$ws1.Range("A4").Value = "ADANIPORTS"
$ws1.Range("B4").Value = 1479.5
$ws1.Range("C4").Value = 1489
$ws1.Range("D4").Value = 1462.5
$ws1.Range("E4").Value = 1478
$ws1.Range("F4").Value = 3290072
$ws1.Range("G4").Value = 6512312
$ws1.Range("H4").Value = -0.4947920185642211
$ws1.Range("I4").Value = "ADANIPORTS"

$ws1.Range("A5").Value = "KOTAKBANK"
$ws1.Range("B5").Value = 2203.9
$ws1.Range("C5").Value = 2209.9
$ws1.Range("D5").Value = 2190.4
$ws1.Range("E5").Value = 2200
$ws1.Range("F5").Value = 3341302
$ws1.Range("G5").Value = 6738420
$ws1.Range("H5").Value = -0.5041416237040731
$ws1.Range("I5").Value = "KOTAKBANK"

$ws1.Range("A6").Value = "BAJAJFINSV"
$ws1.Range("B6").Value = 2090
$ws1.Range("C6").Value = 2104
$ws1.Range("D6").Value = 2077
$ws1.Range("E6").Value = 2081.1
$ws1.Range("F6").Value = 609393
$ws1.Range("G6").Value = 1424932
$ws1.Range("H6").Value = -0.5723353816182105
$ws1.Range("I6").Value = "BAJAJFINSV"

$ws1.Range("A7").Value = "ICICIGI"
$ws1.Range("B7").Value = 1996.5
$ws1.Range("C7").Value = 2017.5
$ws1.Range("D7").Value = 1991.6
$ws1.Range("E7").Value = 2006
$ws1.Range("F7").Value = 355745
$ws1.Range("G7").Value = 860007
$ws1.Range("H7").Value = -0.5863463902038006
$ws1.Range("I7").Value = "ICICIGI"

$ws1.Range("A8").Value = "BSE"
$ws1.Range("B8").Value = 2512
$ws1.Range("C8").Value = 2538
$ws1.Range("D8").Value = 2461.8
$ws1.Range("E8").Value = 2488.2
$ws1.Range("F8").Value = 3965773
$ws1.Range("G8").Value = 9183672
$ws1.Range("H8").Value = -0.5681713153518549
$ws1.Range("I8").Value = "BSE"

$ws1.Range("A9").Value = "SBICARD"
$ws1.Range("B9").Value = 931.05
$ws1.Range("C9").Value = 939.6
$ws1.Range("D9").Value = 923.1
$ws1.Range("E9").Value = 926.7
$ws1.Range("F9").Value = 292090
$ws1.Range("G9").Value = 677803
$ws1.Range("H9").Value = -0.5690635774701499
$ws1.Range("I9").Value = "SBICARD"

$ws1.Range("A10").Value = "PERSISTENT"
$ws1.Range("B10").Value = 5800
$ws1.Range("C10").Value = 5825
$ws1.Range("D10").Value = 5725
$ws1.Range("E10").Value = 5758.1
$ws1.Range("F10").Value = 469353
$ws1.Range("G10").Value = 940813
$ws1.Range("H10").Value = -0.5011197761935687
$ws1.Range("I10").Value = "PERSISTENT"

$ws1.Range("A11").Value = "ANGELONE"
$ws1.Range("B11").Value = 2468.1
$ws1.Range("C11").Value = 2559
$ws1.Range("D11").Value = 2460.7
$ws1.Range("E11").Value = 2495
$ws1.Range("F11").Value = 1585685
$ws1.Range("G11").Value = 3242812
$ws1.Range("H11").Value = -0.5110154396862969
$ws1.Range("I11").Value = "ANGELONE"

$ws1.Range("A12").Value = "CESC"
$ws1.Range("B12").Value = 174.76
$ws1.Range("C12").Value = 174.76
$ws1.Range("D12").Value = 170
$ws1.Range("E12").Value = 170.4
$ws1.Range("F12").Value = 1819909
$ws1.Range("G12").Value = 3775700
$ws1.Range("H12").Value = -0.5179942792065048
$ws1.Range("I12").Value = "CESC"

$ws2 = $wb.Worksheets.Item("Pos_Change")

$ws2.Range("A2").Value = "M&M"
$ws2.Range("B2").Value = 3561
$ws2.Range("C2").Value = 3656
$ws2.Range("D2").Value = 3560.1
$ws2.Range("E2").Value = 3656
$ws2.Range("F2").Value = 3339579
$ws2.Range("G2").Value = 2125552
$ws2.Range("H2").Value = 0.5711584567208895
$ws2.Range("I2").Value = "M&M"

$ws2.Range("A3").Value = "HINDUNILVR"
$ws2.Range("B3").Value = 2565.5
$ws2.Range("C3").Value = 2615
$ws2.Range("D3").Value = 2560
$ws2.Range("E3").Value = 2603.4
$ws2.Range("F3").Value = 2050690
$ws2.Range("G3").Value = 1446646
$ws2.Range("H3").Value = 0.4175479004538775
$ws2.Range("I3").Value = "HINDUNILVR"

$ws2.Range("A4").Value = "RELIANCE"
$ws2.Range("B4").Value = 1401
$ws2.Range("C4").Value = 1423.3
$ws2.Range("D4").Value = 1399.1
$ws2.Range("E4").Value = 1419.1
$ws2.Range("F4").Value = 19335561
$ws2.Range("G4").Value = 12315932
$ws2.Range("H4").Value = 0.5699632800830664
$ws2.Range("I4").Value = "RELIANCE"

$ws2.Range("A5").Value = "MARUTI"
$ws2.Range("B5").Value = 16298
$ws2.Range("C5").Value = 16549
$ws2.Range("D5").Value = 16298
$ws2.Range("E5").Value = 16380
$ws2.Range("F5").Value = 411608
$ws2.Range("G5").Value = 265727
$ws2.Range("H5").Value = 0.5489882473365522
$ws2.Range("I5").Value = "MARUTI"

$ws2.Range("A6").Value = "ULTRACEMCO"
$ws2.Range("B6").Value = 12300
$ws2.Range("C6").Value = 12390
$ws2.Range("D6").Value = 12259
$ws2.Range("E6").Value = 12313
$ws2.Range("F6").Value = 182242
$ws2.Range("G6").Value = 129372
$ws2.Range("H6").Value = 0.4086664811551186
$ws2.Range("I6").Value = "ULTRACEMCO"

$ws2.Range("A7").Value = "HCLTECH"
$ws2.Range("B7").Value = 1499
$ws2.Range("C7").Value = 1501.3
$ws2.Range("D7").Value = 1479.2
$ws2.Range("E7").Value = 1484.9
$ws2.Range("F7").Value = 3210161
$ws2.Range("G7").Value = 2021043
$ws2.Range("H7").Value = 0.5883684810268758
$ws2.Range("I7").Value = "HCLTECH"

$ws2.Range("A8").Value = "TVSMOTOR"
$ws2.Range("B8").Value = 3576.5
$ws2.Range("C8").Value = 3658
$ws2.Range("D8").Value = 3575.5
$ws2.Range("E8").Value = 3653.2
$ws2.Range("F8").Value = 1173044
$ws2.Range("G8").Value = 757797
$ws2.Range("H8").Value = 0.5479660120058538
$ws2.Range("I8").Value = "TVSMOTOR"

$ws2.Range("A9").Value = "NAUKRI"
$ws2.Range("B9").Value = 1345.2
$ws2.Range("C9").Value = 1347
$ws2.Range("D9").Value = 1323
$ws2.Range("E9").Value = 1330.1
$ws2.Range("F9").Value = 904220
$ws2.Range("G9").Value = 595870
$ws2.Range("H9").Value = 0.5174786446708175
$ws2.Range("I9").Value = "NAUKRI"

$ws2.Range("A10").Value = "INDUSINDBK"
$ws2.Range("B10").Value = 740
$ws2.Range("C10").Value = 761
$ws2.Range("D10").Value = 734.5
$ws2.Range("E10").Value = 752.8
$ws2.Range("F10").Value = 4770012
$ws2.Range("G10").Value = 3056360
$ws2.Range("H10").Value = 0.560683950843487
$ws2.Range("I10").Value = "INDUSINDBK"

$ws2.Range("A11").Value = "MARICO"
$ws2.Range("B11").Value = 728.5
$ws2.Range("C11").Value = 739.5
$ws2.Range("D11").Value = 726.95
$ws2.Range("E11").Value = 736
$ws2.Range("F11").Value = 1707666
$ws2.Range("G11").Value = 1205842
$ws2.Range("H11").Value = 0.4161606578639656
$ws2.Range("I11").Value = "MARICO"

$ws2.Range("A12").Value = "KALYANKJIL"
$ws2.Range("B12").Value = 487.05
$ws2.Range("C12").Value = 502.5
$ws2.Range("D12").Value = 485.85
$ws2.Range("E12").Value = 490.25
$ws2.Range("F12").Value = 6416427
$ws2.Range("G12").Value = 4326210
$ws2.Range("H12").Value = 0.4831519967824031
$ws2.Range("I12").Value = "KALYANKJIL"

$ws2.Range("A13").Value = "HEROMOTOCO"
$ws2.Range("B13").Value = 5582
$ws2.Range("C13").Value = 5625
$ws2.Range("D13").Value = 5550.5
$ws2.Range("E13").Value = 5595
$ws2.Range("F13").Value = 435204
$ws2.Range("G13").Value = 294066
$ws2.Range("H13").Value = 0.4799534798310583
$ws2.Range("I13").Value = "HEROMOTOCO"

$ws2.Range("A14").Value = "PAGEIND"
$ws2.Range("B14").Value = 41270
$ws2.Range("C14").Value = 41270
$ws2.Range("D14").Value = 40780
$ws2.Range("E14").Value = 41050
$ws2.Range("F14").Value = 21226
$ws2.Range("G14").Value = 14326
$ws2.Range("H14").Value = 0.4816417702080134
$ws2.Range("I14").Value = "PAGEIND"

$ws2.Range("A15").Value = "SUPREMEIND"
$ws2.Range("B15").Value = 4240.2
$ws2.Range("C15").Value = 4266
$ws2.Range("D15").Value = 4179
$ws2.Range("E15").Value = 4220
$ws2.Range("F15").Value = 43475
$ws2.Range("G15").Value = 29283
$ws2.Range("H15").Value = 0.4846497968104361
$ws2.Range("I15").Value = "SUPREMEIND"

$ws2.Range("A16").Value = "TATAELXSI"
$ws2.Range("B16").Value = 5403.5
$ws2.Range("C16").Value = 5420
$ws2.Range("D16").Value = 5335.5
$ws2.Range("E16").Value = 5373.5
$ws2.Range("F16").Value = 193891
$ws2.Range("G16").Value = 122296
$ws2.Range("H16").Value = 0.5854238895793812
$ws2.Range("I16").Value = "TATAELXSI"

$ws2.Range("A17").Value = "UNIONBANK"
$ws2.Range("B17").Value = 138.81
$ws2.Range("C17").Value = 139.08
$ws2.Range("D17").Value = 135.8
$ws2.Range("E17").Value = 136.77
$ws2.Range("F17").Value = 7584213
$ws2.Range("G17").Value = 5116871
$ws2.Range("H17").Value = 0.4821974210411011
$ws2.Range("I17").Value = "UNIONBANK"

$ws2.Range("A18").Value = "TIINDIA"
$ws2.Range("B18").Value = 3179
$ws2.Range("C18").Value = 3189.4
$ws2.Range("D18").Value = 3108.9
$ws2.Range("E18").Value = 3130
$ws2.Range("F18").Value = 180720
$ws2.Range("G18").Value = 127384
$ws2.Range("H18").Value = 0.4187025058092068
$ws2.Range("I18").Value = "TIINDIA"

$ws2.Range("A19").Value = "ASHOKLEY"
$ws2.Range("B19").Value = 137.35
$ws2.Range("C19").Value = 137.7
$ws2.Range("D19").Value = 134.21
$ws2.Range("E19").Value = 134.52
$ws2.Range("F19").Value = 15168999
$ws2.Range("G19").Value = 9950107
$ws2.Range("H19").Value = 0.5245061183764155
$ws2.Range("I19").Value = "ASHOKLEY"

$ws2.Range("A20").Value = "POLICYBZR"
$ws2.Range("B20").Value = 1689.1
$ws2.Range("C20").Value = 1689.1
$ws2.Range("D20").Value = 1637.1
$ws2.Range("E20").Value = 1647.6
$ws2.Range("F20").Value = 2926578
$ws2.Range("G20").Value = 1990102
$ws2.Range("H20").Value = 0.4705668352677401
$ws2.Range("I20").Value = "POLICYBZR"

$ws2.Range("A21").Value = "GRANULES"
$ws2.Range("B21").Value = 570.5
$ws2.Range("C21").Value = 574.7
$ws2.Range("D21").Value = 567.05
$ws2.Range("E21").Value = 567.5
$ws2.Range("F21").Value = 727221
$ws2.Range("G21").Value = 508720
$ws2.Range("H21").Value = 0.4295113225349898
$ws2.Range("I21").Value = "GRANULES"

